# Refresh scraped schedule data for Linea 141 (scrape run @ 06:25:43)
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 06:25:43"
$ws.Range("A3").Value = "Total filas: 56"

$ws.Cells.Item(26, 1).Value = "06:25:43"
$ws.Cells.Item(26, 2).Value = "06:26"
$ws.Cells.Item(26, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = "LP1912"
$ws.Cells.Item(27, 1).Value = "04:45:05"
$ws.Cells.Item(27, 2).Value = "06:27"
$ws.Cells.Item(27, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(27, 4).Value = 102
$ws.Cells.Item(27, 5).Value = "LP1912"
$ws.Cells.Item(28, 1).Value = "06:25:43"
$ws.Cells.Item(28, 2).Value = "06:28"
$ws.Cells.Item(28, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 5).Value = "LP1912"
$ws.Cells.Item(29, 1).Value = "04:56:49"
$ws.Cells.Item(29, 2).Value = "06:29"
$ws.Cells.Item(29, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(29, 4).Value = 93
$ws.Cells.Item(29, 5).Value = "LP1912"
$ws.Cells.Item(30, 1).Value = "04:45:05"
$ws.Cells.Item(30, 2).Value = "06:30"
$ws.Cells.Item(30, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(30, 4).Value = 105
$ws.Cells.Item(30, 5).Value = "LP1912"
$ws.Cells.Item(31, 1).Value = "04:45:05"
$ws.Cells.Item(31, 2).Value = "06:31"
$ws.Cells.Item(31, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(31, 4).Value = 106
$ws.Cells.Item(31, 5).Value = "LP1912"
$ws.Cells.Item(32, 1).Value = "04:45:05"
$ws.Cells.Item(32, 2).Value = "06:44"
$ws.Cells.Item(32, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(32, 4).Value = 119
$ws.Cells.Item(32, 5).Value = "LP1912"
$ws.Cells.Item(33, 1).Value = "05:55:25"
$ws.Cells.Item(33, 2).Value = "06:44"
$ws.Cells.Item(33, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(33, 4).Value = 49
$ws.Cells.Item(33, 5).Value = "LP1912"
$ws.Cells.Item(34, 1).Value = "04:56:49"
$ws.Cells.Item(34, 2).Value = "06:46"
$ws.Cells.Item(34, 3).Value = "215C_EL PATO"
$ws.Cells.Item(34, 4).Value = 110
$ws.Cells.Item(34, 5).Value = "LP1912"
$ws.Cells.Item(35, 1).Value = "05:26:08"
$ws.Cells.Item(35, 2).Value = "06:47"
$ws.Cells.Item(35, 3).Value = "215C_EL PATO"
$ws.Cells.Item(35, 4).Value = 81
$ws.Cells.Item(35, 5).Value = "LP1912"
$ws.Cells.Item(36, 1).Value = "05:55:25"
$ws.Cells.Item(36, 2).Value = "06:59"
$ws.Cells.Item(36, 3).Value = "14_ABASTO"
$ws.Cells.Item(36, 4).Value = 64
$ws.Cells.Item(36, 5).Value = "LP1912"
$ws.Cells.Item(37, 1).Value = "05:26:08"
$ws.Cells.Item(37, 2).Value = "07:00"
$ws.Cells.Item(37, 3).Value = "14_ABASTO"
$ws.Cells.Item(37, 4).Value = 94
$ws.Cells.Item(37, 5).Value = "LP1912"
$ws.Cells.Item(38, 1).Value = "06:25:43"
$ws.Cells.Item(38, 2).Value = "07:01"
$ws.Cells.Item(38, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(38, 4).Value = 36
$ws.Cells.Item(38, 5).Value = "LP1912"
$ws.Cells.Item(39, 1).Value = "05:55:25"
$ws.Cells.Item(39, 2).Value = "07:04"
$ws.Cells.Item(39, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(39, 4).Value = 69
$ws.Cells.Item(39, 5).Value = "LP1912"
$ws.Cells.Item(40, 1).Value = "05:26:08"
$ws.Cells.Item(40, 2).Value = "07:05"
$ws.Cells.Item(40, 3).Value = "15_ABASTO"
$ws.Cells.Item(40, 4).Value = 99
$ws.Cells.Item(40, 5).Value = "LP1912"
$ws.Cells.Item(41, 1).Value = "05:26:08"
$ws.Cells.Item(41, 2).Value = "07:05"
$ws.Cells.Item(41, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(41, 4).Value = 99
$ws.Cells.Item(41, 5).Value = "LP1912"
$ws.Cells.Item(42, 1).Value = "05:26:08"
$ws.Cells.Item(42, 2).Value = "07:06"
$ws.Cells.Item(42, 3).Value = "10_OLMOS"
$ws.Cells.Item(42, 4).Value = 100
$ws.Cells.Item(42, 5).Value = "LP1912"
$ws.Cells.Item(43, 1).Value = "05:26:08"
$ws.Cells.Item(43, 2).Value = "07:07"
$ws.Cells.Item(43, 3).Value = "225_GOMEZ"
$ws.Cells.Item(43, 4).Value = 101
$ws.Cells.Item(43, 5).Value = "LP1912"
$ws.Cells.Item(44, 1).Value = "05:26:08"
$ws.Cells.Item(44, 2).Value = "07:11"
$ws.Cells.Item(44, 3).Value = "215A_EL PATO"
$ws.Cells.Item(44, 4).Value = 105
$ws.Cells.Item(44, 5).Value = "LP1912"
$ws.Cells.Item(45, 1).Value = "06:25:43"
$ws.Cells.Item(45, 2).Value = "07:14"
$ws.Cells.Item(45, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(45, 4).Value = 49
$ws.Cells.Item(45, 5).Value = "LP1912"
$ws.Cells.Item(46, 1).Value = "05:26:08"
$ws.Cells.Item(46, 2).Value = "07:15"
$ws.Cells.Item(46, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(46, 4).Value = 80
$ws.Cells.Item(46, 5).Value = "LP1912"
$ws.Cells.Item(47, 1).Value = "05:26:08"
$ws.Cells.Item(47, 2).Value = "07:16"
$ws.Cells.Item(47, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(47, 4).Value = 110
$ws.Cells.Item(47, 5).Value = "LP1912"
$ws.Cells.Item(48, 1).Value = "05:26:08"
$ws.Cells.Item(48, 2).Value = "07:21"
$ws.Cells.Item(48, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(48, 4).Value = 115
$ws.Cells.Item(48, 5).Value = "LP1912"
$ws.Cells.Item(49, 1).Value = "05:26:08"
$ws.Cells.Item(49, 2).Value = "07:23"
$ws.Cells.Item(49, 3).Value = "10_OLMOS"
$ws.Cells.Item(49, 4).Value = 117
$ws.Cells.Item(49, 5).Value = "LP1912"
$ws.Cells.Item(50, 1).Value = "05:55:25"
$ws.Cells.Item(50, 2).Value = "07:30"
$ws.Cells.Item(50, 3).Value = "10_OLMOS"
$ws.Cells.Item(50, 4).Value = 95
$ws.Cells.Item(50, 5).Value = "LP1912"
$ws.Cells.Item(51, 1).Value = "05:55:25"
$ws.Cells.Item(51, 2).Value = "07:31"
$ws.Cells.Item(51, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(51, 4).Value = 96
$ws.Cells.Item(51, 5).Value = "LP1912"
$ws.Cells.Item(52, 1).Value = "05:55:25"
$ws.Cells.Item(52, 2).Value = "07:31"
$ws.Cells.Item(52, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(52, 4).Value = 96
$ws.Cells.Item(52, 5).Value = "LP1912"
$ws.Cells.Item(53, 1).Value = "05:55:25"
$ws.Cells.Item(53, 2).Value = "07:32"
$ws.Cells.Item(53, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(53, 4).Value = 97
$ws.Cells.Item(53, 5).Value = "LP1912"
$ws.Cells.Item(54, 1).Value = "05:55:25"
$ws.Cells.Item(54, 2).Value = "07:36"
$ws.Cells.Item(54, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(54, 4).Value = 101
$ws.Cells.Item(54, 5).Value = "LP1912"
$ws.Cells.Item(55, 1).Value = "05:55:25"
$ws.Cells.Item(55, 2).Value = "07:39"
$ws.Cells.Item(55, 3).Value = "10_OLMOS"
$ws.Cells.Item(55, 4).Value = 104
$ws.Cells.Item(55, 5).Value = "LP1912"
$ws.Cells.Item(56, 1).Value = "05:55:25"
$ws.Cells.Item(56, 2).Value = "07:47"
$ws.Cells.Item(56, 3).Value = "14_ABASTO"
$ws.Cells.Item(56, 4).Value = 112
$ws.Cells.Item(56, 5).Value = "LP1912"
$ws.Cells.Item(57, 1).Value = "05:55:25"
$ws.Cells.Item(57, 2).Value = "07:51"
$ws.Cells.Item(57, 3).Value = "215D_EL PATO"
$ws.Cells.Item(57, 4).Value = 116
$ws.Cells.Item(57, 5).Value = "LP1912"
$ws.Cells.Item(58, 1).Value = "06:25:43"
$ws.Cells.Item(58, 2).Value = "08:01"
$ws.Cells.Item(58, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(58, 4).Value = 96
$ws.Cells.Item(58, 5).Value = "LP1912"
$ws.Cells.Item(59, 1).Value = "06:25:43"
$ws.Cells.Item(59, 2).Value = "08:12"
$ws.Cells.Item(59, 3).Value = "15_ABASTO"
$ws.Cells.Item(59, 4).Value = 107
$ws.Cells.Item(59, 5).Value = "LP1912"
$ws.Cells.Item(60, 1).Value = "06:25:43"
$ws.Cells.Item(60, 2).Value = "08:22"
$ws.Cells.Item(60, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(60, 4).Value = 117
$ws.Cells.Item(60, 5).Value = "LP1912"
$ws.Cells.Item(61, 1).Value = "06:25:43"
$ws.Cells.Item(61, 2).Value = "08:23"
$ws.Cells.Item(61, 3).Value = "215B_EL PATO"
$ws.Cells.Item(61, 4).Value = 118
$ws.Cells.Item(61, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 06:25:43"
$ws.Range("A3").Value = "Total filas: 10"

$ws.Cells.Item(15, 1).Value = "06:25:43"
$ws.Cells.Item(15, 2).Value = "08:23"
$ws.Cells.Item(15, 3).Value = "215B_EL PATO"
$ws.Cells.Item(15, 4).Value = 118
$ws.Cells.Item(15, 5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 06:25:43"
$ws.Range("A3").Value = "Total filas: 11"

$ws.Cells.Item(12, 1).Value = "06:25:43"
$ws.Cells.Item(12, 2).Value = "06:59"
$ws.Cells.Item(12, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(12, 4).Value = 34
$ws.Cells.Item(12, 5).Value = "L6173"
$ws.Cells.Item(13, 1).Value = "05:26:08"
$ws.Cells.Item(13, 2).Value = "07:00"
$ws.Cells.Item(13, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(13, 4).Value = 94
$ws.Cells.Item(13, 5).Value = "L6173"
$ws.Cells.Item(14, 1).Value = "05:55:25"
$ws.Cells.Item(14, 2).Value = "07:35"
$ws.Cells.Item(14, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(14, 4).Value = 100
$ws.Cells.Item(14, 5).Value = "L6173"
$ws.Cells.Item(15, 1).Value = "06:25:43"
$ws.Cells.Item(15, 2).Value = "07:39"
$ws.Cells.Item(15, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(15, 4).Value = 74
$ws.Cells.Item(15, 5).Value = "L6173"
$ws.Cells.Item(16, 1).Value = "06:25:43"
$ws.Cells.Item(16, 2).Value = "08:06"
$ws.Cells.Item(16, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(16, 4).Value = 101
$ws.Cells.Item(16, 5).Value = "L6203"

